$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resize")
$ws.Range("AK1").Value = "node 234mb zip"
$ws.Range("AK1").Interior.Color = 255
$ws.Range("AK1").Font.Color = 255
